$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 209, shifting existing rows 209:326 down to 210:327.
$ws.Rows.Item(209).Insert()

# Populate the newly inserted row 209 with the new record's data.
$ws.Cells.Item(209, 1).Value = 9
$ws.Cells.Item(209, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(209, 3).Value = "Metropolitana"
$ws.Cells.Item(209, 4).Value = 44572
$ws.Cells.Item(209, 5).Value = 13
$ws.Cells.Item(209, 6).Value = 100112039
$ws.Cells.Item(209, 7).Value = "Ciboulette"
$ws.Cells.Item(209, 8).Value = "Sin especificar"
$ws.Cells.Item(209, 9).Value = "Primera"
$ws.Cells.Item(209, 10).Value = 250
$ws.Cells.Item(209, 11).Value = 1000
$ws.Cells.Item(209, 12).Value = 1200
$ws.Cells.Item(209, 13).Value = 1100
$ws.Cells.Item(209, 14).Value = "`$/docena de atados"
$ws.Cells.Item(209, 15).Value = "Región Metropolitana"
$ws.Cells.Item(209, 16).Value = 367
$ws.Cells.Item(209, 17).Value = 3
$ws.Cells.Item(209, 18).Value = "Hortaliza"
